$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: set C1 value and copy style/format from B1
$ws.Cells.Item(1, 3).Value = "13-01-2023"
$ws.Range("B1").Copy() | Out-Null
$ws.Cells.Item(1, 3).PasteSpecial(-4122) | Out-Null

# Rows 2-50: update A and B values (existing styled cells), set new C values
$ws.Cells.Item(2, 1).Value = "1810 Renta variable"
$ws.Cells.Item(2, 2).Value = 104034.32
$ws.Cells.Item(2, 3).Value = 116303.33
$ws.Cells.Item(3, 1).Value = "1822 Raices Valores Negociables"
$ws.Cells.Item(3, 2).Value = 191911.68
$ws.Cells.Item(3, 3).Value = 203468.12
$ws.Cells.Item(4, 1).Value = "Adcap IOL Acciones Argentina"
$ws.Cells.Item(4, 2).Value = 31947.73
$ws.Cells.Item(4, 3).Value = 33986.48
$ws.Cells.Item(5, 1).Value = "Allaria Acciones"
$ws.Cells.Item(5, 2).Value = 53609.39
$ws.Cells.Item(5, 3).Value = 53649.92
$ws.Cells.Item(6, 1).Value = "Alpha Acciones"
$ws.Cells.Item(6, 2).Value = 84343.23
$ws.Cells.Item(6, 3).Value = 84284.61
$ws.Cells.Item(7, 1).Value = "Alpha Mega"
$ws.Cells.Item(7, 2).Value = 232798.17
$ws.Cells.Item(7, 3).Value = 232958.18
$ws.Cells.Item(8, 1).Value = "Alpha planeam equil"
$ws.Cells.Item(8, 2).Value = 6662.21
$ws.Cells.Item(8, 3).Value = 4153.07
$ws.Cells.Item(9, 1).Value = "Alpha renta balan global"
$ws.Cells.Item(9, 2).Value = 704406.59
$ws.Cells.Item(9, 3).Value = 701809.37
$ws.Cells.Item(10, 1).Value = "Argenfunds"
$ws.Cells.Item(10, 2).Value = 10072.37
$ws.Cells.Item(10, 3).Value = 10064.37
$ws.Cells.Item(11, 1).Value = "Arpenta ex Mercosur"
$ws.Cells.Item(11, 2).Value = 16812.59
$ws.Cells.Item(11, 3).Value = 16817.52
$ws.Cells.Item(12, 1).Value = "Balanz"
$ws.Cells.Item(12, 2).Value = 208523.69
$ws.Cells.Item(12, 3).Value = 176487.81
$ws.Cells.Item(13, 1).Value = "Bull Market"
$ws.Cells.Item(13, 2).Value = 59095.22
$ws.Cells.Item(13, 3).Value = 66018.27
$ws.Cells.Item(14, 1).Value = "CMA acciones"
$ws.Cells.Item(14, 2).Value = 213757.81
$ws.Cells.Item(14, 3).Value = 213795.36
$ws.Cells.Item(15, 1).Value = "Compass Crecimiento"
$ws.Cells.Item(15, 2).Value = 1284814.59
$ws.Cells.Item(15, 3).Value = 1272281.84
$ws.Cells.Item(16, 1).Value = "Compass Crecimiento II"
$ws.Cells.Item(16, 2).Value = 12520.43
$ws.Cells.Item(16, 3).Value = 12520.97
$ws.Cells.Item(17, 1).Value = "Consultatio Acciones Argentina"
$ws.Cells.Item(17, 2).Value = 439565.08
$ws.Cells.Item(17, 3).Value = 440084.49
$ws.Cells.Item(18, 1).Value = "Consultatio Renta Variable"
$ws.Cells.Item(18, 2).Value = 174898.74
$ws.Cells.Item(18, 3).Value = 175142.86
$ws.Cells.Item(19, 1).Value = "Delta Acciones"
$ws.Cells.Item(19, 2).Value = 61968.7
$ws.Cells.Item(19, 3).Value = 61998.59
$ws.Cells.Item(20, 1).Value = "Delta Internacional"
$ws.Cells.Item(20, 2).Value = 168.74
$ws.Cells.Item(20, 3).Value = 167.24
$ws.Cells.Item(21, 1).Value = "Delta Recursos Naturales"
$ws.Cells.Item(21, 2).Value = 558129.01
$ws.Cells.Item(21, 3).Value = 557863.1899999999
$ws.Cells.Item(22, 1).Value = "Delta Select"
$ws.Cells.Item(22, 2).Value = 390493.53
$ws.Cells.Item(22, 3).Value = 390715.62
$ws.Cells.Item(23, 1).Value = "Delta gestion V"
$ws.Cells.Item(23, 2).Value = 106145.82
$ws.Cells.Item(23, 3).Value = 106308.8
$ws.Cells.Item(24, 1).Value = "FBA Acciones Argentinas"
$ws.Cells.Item(24, 2).Value = 175952.1
$ws.Cells.Item(24, 3).Value = 179588.92
$ws.Cells.Item(25, 1).Value = "FBA Calificado"
$ws.Cells.Item(25, 2).Value = 173673.42
$ws.Cells.Item(25, 3).Value = 175981.38
$ws.Cells.Item(26, 1).Value = "Fima Acciones"
$ws.Cells.Item(26, 2).Value = 215351.45
$ws.Cells.Item(26, 3).Value = 232196.75
$ws.Cells.Item(27, 1).Value = "Fima PB Acciones"
$ws.Cells.Item(27, 2).Value = 95490.92999999999
$ws.Cells.Item(27, 3).Value = 100237.12
$ws.Cells.Item(28, 1).Value = "Gainvest Renta Variable"
$ws.Cells.Item(28, 2).Value = 87729.00999999999
$ws.Cells.Item(28, 3).Value = 87721.81
$ws.Cells.Item(29, 1).Value = "Galileo Acciones"
$ws.Cells.Item(29, 2).Value = 1976674.58
$ws.Cells.Item(29, 3).Value = 1928360.97
$ws.Cells.Item(30, 1).Value = "Goal Acciones Argentinas"
$ws.Cells.Item(30, 2).Value = 45747.36
$ws.Cells.Item(30, 3).Value = 45714.74
$ws.Cells.Item(31, 1).Value = "Goal acciones plus"
$ws.Cells.Item(31, 2).Value = 4558.05
$ws.Cells.Item(31, 3).Value = 4555.44
$ws.Cells.Item(32, 1).Value = "HF Acciones Argentinas"
$ws.Cells.Item(32, 2).Value = 138789.82
$ws.Cells.Item(32, 3).Value = 133733.32
$ws.Cells.Item(33, 1).Value = "HF Acciones Lideres"
$ws.Cells.Item(33, 2).Value = 275283.19
$ws.Cells.Item(33, 3).Value = 265073.7
$ws.Cells.Item(34, 1).Value = "IAM Renta Variable"
$ws.Cells.Item(34, 2).Value = 32436.65
$ws.Cells.Item(34, 3).Value = 33943.21
$ws.Cells.Item(35, 1).Value = "IEB Value"
$ws.Cells.Item(35, 2).Value = 8103.05
$ws.Cells.Item(35, 3).Value = 8104.65
$ws.Cells.Item(36, 1).Value = "Lombardi"
$ws.Cells.Item(36, 2).Value = 21706.76
$ws.Cells.Item(36, 3).Value = 21711.31
$ws.Cells.Item(37, 1).Value = "MAF"
$ws.Cells.Item(37, 2).Value = 23745.09
$ws.Cells.Item(37, 3).Value = 23704.47
$ws.Cells.Item(38, 1).Value = "Megainver"
$ws.Cells.Item(38, 2).Value = 30176.82
$ws.Cells.Item(38, 3).Value = 30163.52
$ws.Cells.Item(39, 1).Value = "Pellegrini Acciones"
$ws.Cells.Item(39, 2).Value = 61581.46
$ws.Cells.Item(39, 3).Value = 61628.97
$ws.Cells.Item(40, 1).Value = "Pionero Acciones"
$ws.Cells.Item(40, 2).Value = 132762.87
$ws.Cells.Item(40, 3).Value = 132740.58
$ws.Cells.Item(41, 1).Value = "Premier Renta Variable"
$ws.Cells.Item(41, 2).Value = 34754.06
$ws.Cells.Item(41, 3).Value = 34785.32
$ws.Cells.Item(42, 1).Value = "Quinquela Acciones"
$ws.Cells.Item(42, 2).Value = 96329.37
$ws.Cells.Item(42, 3).Value = 96227.52
$ws.Cells.Item(43, 1).Value = "Rofex 20 Renta Variable"
$ws.Cells.Item(43, 2).Value = 66545.2
$ws.Cells.Item(43, 3).Value = 66517.99000000001
$ws.Cells.Item(44, 1).Value = "SBS Acciones Argentina"
$ws.Cells.Item(44, 2).Value = 548193.05
$ws.Cells.Item(44, 3).Value = 560367.16
$ws.Cells.Item(45, 1).Value = "Schroeder RV"
$ws.Cells.Item(45, 2).Value = 1317408.94
$ws.Cells.Item(45, 3).Value = 1317359.35
$ws.Cells.Item(46, 1).Value = "Supefondo RV"
$ws.Cells.Item(46, 2).Value = 715172.0600000001
$ws.Cells.Item(46, 3).Value = 716517.87
$ws.Cells.Item(47, 1).Value = "Superfondo "
$ws.Cells.Item(47, 2).Value = 54956.5
$ws.Cells.Item(47, 3).Value = 55248.68
$ws.Cells.Item(48, 1).Value = "Toronto Trust Multimercado"
$ws.Cells.Item(48, 2).Value = 41918.65
$ws.Cells.Item(48, 3).Value = 38263.48
$ws.Cells.Item(49, 1).Value = "avg"
$ws.Cells.Item(49, 2).Value = 240887.66
$ws.Cells.Item(49, 3).Value = 240028.26
$ws.Cells.Item(50, 1).Value = "total"
$ws.Cells.Item(50, 2).Value = 11321720.08
$ws.Cells.Item(50, 3).Value = 11281328.24

Write-Host "done"